# Update pembayaran Purwadi & Ipin
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Row 105: add Tips (L) value ---
$ws.Range("L105").Value = 150000

# --- Expand the table by two rows (106 already existed inside the table range; ---
# --- we need two brand-new rows, 107 and 108) ---
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# --- Row 106: Purwadi payment ---
$ws.Range("A106").Value = 45387
$ws.Range("B106").Value = "Purwadi"
$ws.Range("D106").Value = 45373
$ws.Range("E106").Value = 45386
$ws.Range("F106").Value = 3718000
$ws.Range("G106").Value = 3718000
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 3657000
$ws.Range("K106").Formula = "=G106+H106-J106"
$ws.Range("L106").Value = 100000

# --- Row 107: Perorangan-2 / Ipin payment ---
$ws.Range("A107").NumberFormat = "dd/mm/yyyy"
$ws.Range("B107").Value = "Perorangan-2"
$ws.Range("C107").Value = "Ipin"
$ws.Range("D107").NumberFormat = "dd/mm/yyyy"
$ws.Range("D107").Value = 45381
$ws.Range("E107").NumberFormat = "dd/mm/yyyy"
$ws.Range("E107").Value = 45386
$ws.Range("F107").Value = 225000
$ws.Range("G107").Value = 225000
$ws.Range("H107").Value = 0
$ws.Range("I107").Formula = "=F107-G107"
$ws.Range("J107").Value = 225000
$ws.Range("K107").Formula = "=G107+H107-J107"

# --- Row 108: blank placeholder row (kept for formatting/formula continuity) ---
$ws.Range("A108").NumberFormat = "dd/mm/yyyy"
$ws.Range("D108").NumberFormat = "dd/mm/yyyy"
$ws.Range("E108").NumberFormat = "dd/mm/yyyy"
$ws.Range("I108").Formula = "=F108-G108"

# --- Update view state to match where the user left off editing ---
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G101").Select()

$wb.Save()
